# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple Leve-profit worksheets, refreshing cached price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 397.92307
$ws.Range("I53").Value = 324.875
$ws.Range("K53").Value = 324.875
$ws.Range("M53").Value = 312.125
$ws.Range("H106").Value = 5124.25
$ws.Range("I106").Value = 5505
$ws.Range("K106").Value = 5505
$ws.Range("M106").Value = -4874
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H127").Value = 2423.125
$ws.Range("J127").Value = 2349.6667
$ws.Range("L127").Value = 7049.000100000001
$ws.Range("N127").Value = -16969.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4208.3335
$ws.Range("I45").Value = 3944.4443
$ws.Range("K45").Value = 3944.4443
$ws.Range("M45").Value = -3567.4443
$ws.Range("H61").Value = 7599.8
$ws.Range("I61").Value = 1833
$ws.Range("J61").Value = 16250
$ws.Range("K61").Value = 1833
$ws.Range("L61").Value = 16250
$ws.Range("M61").Value = -1621
$ws.Range("N61").Value = -16674
$ws.Range("H74").Value = 902.8
$ws.Range("I74").Value = 893.8
$ws.Range("K74").Value = 893.8
$ws.Range("M74").Value = -19.79999999999995
$ws.Range("H77").Value = 902.8
$ws.Range("I77").Value = 893.8
$ws.Range("K77").Value = 4469
$ws.Range("M77").Value = -101
$ws.Range("H132").Value = 1432
$ws.Range("I132").Value = 1432
$ws.Range("K132").Value = 4296
$ws.Range("M132").Value = -1766
$ws.Range("H136").Value = 7599.8
$ws.Range("I136").Value = 1833
$ws.Range("J136").Value = 16250
$ws.Range("K136").Value = 5499
$ws.Range("L136").Value = 48750
$ws.Range("M136").Value = -2949
$ws.Range("N136").Value = -53850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2662
$ws.Range("I105").Value = 2662
$ws.Range("K105").Value = 2662
$ws.Range("M105").Value = -915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2674.5
$ws.Range("I3").Value = 2566
$ws.Range("K3").Value = 2566
$ws.Range("M3").Value = -2453
$ws.Range("H86").Value = 5243.273
$ws.Range("I86").Value = 5442.6
$ws.Range("K86").Value = 5442.6
$ws.Range("M86").Value = -4319.6
$ws.Range("H89").Value = 5243.273
$ws.Range("I89").Value = 5442.6
$ws.Range("K89").Value = 27213
$ws.Range("M89").Value = -21597
$ws.Range("H107").Value = 375.5
$ws.Range("I107").Value = 293.33334
$ws.Range("J107").Value = 457.66666
$ws.Range("K107").Value = 293.33334
$ws.Range("L107").Value = 457.66666
$ws.Range("M107").Value = 1626.66666
$ws.Range("N107").Value = -4297.66666
$ws.Range("H134").Value = 1415.1538
$ws.Range("I134").Value = 1416.4166
$ws.Range("K134").Value = 4249.2498
$ws.Range("M134").Value = -1714.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1999.4
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1999.4
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5998.200000000001
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6166.200000000001
$ws.Range("H39").Value = 1123.4286
$ws.Range("I39").Value = 772.8
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 2318.4
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -2024.4
$ws.Range("N39").Value = -6588
$ws.Range("H109").Value = 1271.3334
$ws.Range("I109").Value = 1271.3334
$ws.Range("K109").Value = 3814.0002
$ws.Range("M109").Value = -2774.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3463.3333
$ws.Range("I80").Value = 3278.3333
$ws.Range("J80").Value = 3833.3333
$ws.Range("K80").Value = 3278.3333
$ws.Range("L80").Value = 3833.3333
$ws.Range("M80").Value = -2280.3333
$ws.Range("N80").Value = -5829.3333
$ws.Range("H83").Value = 3463.3333
$ws.Range("I83").Value = 3278.3333
$ws.Range("J83").Value = 3833.3333
$ws.Range("K83").Value = 16391.6665
$ws.Range("L83").Value = 19166.6665
$ws.Range("M83").Value = -11399.6665
$ws.Range("N83").Value = -29150.6665
$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5840
$ws.Range("H132").Value = 5859
$ws.Range("I132").Value = 5997.3335
$ws.Range("K132").Value = 17992.0005
$ws.Range("M132").Value = -15462.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 786.6667
$ws.Range("I22").Value = 691.5
$ws.Range("J22").Value = 881.8333
$ws.Range("K22").Value = 691.5
$ws.Range("L22").Value = 881.8333
$ws.Range("M22").Value = -396.5
$ws.Range("N22").Value = -1471.8333
$ws.Range("H27").Value = 786.6667
$ws.Range("I27").Value = 691.5
$ws.Range("J27").Value = 881.8333
$ws.Range("K27").Value = 691.5
$ws.Range("L27").Value = 881.8333
$ws.Range("M27").Value = -584.5
$ws.Range("N27").Value = -1095.8333
$ws.Range("H40").Value = 2737.5557
$ws.Range("I40").Value = 2322.9092
$ws.Range("J40").Value = 3389.1428
$ws.Range("K40").Value = 2322.9092
$ws.Range("L40").Value = 3389.1428
$ws.Range("M40").Value = -2186.9092
$ws.Range("N40").Value = -3661.1428
$ws.Range("H93").Value = 745.3461
$ws.Range("I93").Value = 780.375
$ws.Range("J93").Value = 689.3
$ws.Range("K93").Value = 780.375
$ws.Range("L93").Value = 689.3
$ws.Range("M93").Value = 467.625
$ws.Range("N93").Value = -3185.3
$ws.Range("H122").Value = 6111.4634
$ws.Range("I122").Value = 4754.55
$ws.Range("J122").Value = 7403.7617
$ws.Range("K122").Value = 14263.65
$ws.Range("L122").Value = 22211.2851
$ws.Range("M122").Value = -11813.65
$ws.Range("N122").Value = -27111.2851

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 450000
$ws.Range("J2").Value = 882500
$ws.Range("L2").Value = 882500
$ws.Range("N2").Value = -882724
$ws.Range("H3").Value = 25005000
$ws.Range("I3").Value = 50005000
$ws.Range("K3").Value = 50005000
$ws.Range("M3").Value = -50004886
$ws.Range("H7").Value = 1683.3334
$ws.Range("I7").Value = 275
$ws.Range("K7").Value = 275
$ws.Range("M7").Value = -162
$ws.Range("H81").Value = 8820.556
$ws.Range("I81").Value = 8439.799999999999
$ws.Range("J81").Value = 9296.5
$ws.Range("K81").Value = 16879.6
$ws.Range("L81").Value = 18593
$ws.Range("M81").Value = -15818.6
$ws.Range("N81").Value = -20715
$ws.Range("H84").Value = 8820.556
$ws.Range("I84").Value = 8439.799999999999
$ws.Range("J84").Value = 9296.5
$ws.Range("K84").Value = 84398
$ws.Range("L84").Value = 92965
$ws.Range("M84").Value = -79094
$ws.Range("N84").Value = -103573
$ws.Range("H96").Value = 2933.3
$ws.Range("I96").Value = 2819.8572
$ws.Range("K96").Value = 2819.8572
$ws.Range("M96").Value = -1446.8572
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344
$ws.Range("H132").Value = 992.25
$ws.Range("I132").Value = 991.5454999999999
$ws.Range("K132").Value = 2974.6365
$ws.Range("M132").Value = -444.6364999999996
$ws.Range("H136").Value = 2916.7896
$ws.Range("J136").Value = 1947
$ws.Range("L136").Value = 5841
$ws.Range("N136").Value = -10941
